$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1979522184300341
$ws.Range("C2").Value = 0.5324232081911263
$ws.Range("J2").Value = 0.0273037542662116
$ws.Range("P2").Value = 0.1262798634812287
$ws.Range("S2").Value = 0.1160409556313993
$ws.Range("B3").Value = 0.03144654088050314
$ws.Range("C3").Value = 0.006289308176100629
$ws.Range("J3").Value = 0.01886792452830189
$ws.Range("P3").Value = 0.6918238993710691
$ws.Range("S3").Value = 0.2515723270440252
$ws.Range("J4").Value = 0.06060606060606061
$ws.Range("P4").Value = 0.6666666666666666
$ws.Range("S4").Value = 0.2727272727272727
$ws.Range("B6").Value = 0.08673469387755102
$ws.Range("D6").Value = 0.01530612244897959
$ws.Range("F6").Value = 0.07142857142857142
$ws.Range("J6").Value = 0.173469387755102
$ws.Range("O6").Value = 0.01530612244897959
$ws.Range("Q6").Value = 0.1224489795918367
$ws.Range("R6").Value = 0.1071428571428571
$ws.Range("S6").Value = 0.4081632653061225
$ws.Range("B7").Value = 0.1164383561643836
$ws.Range("E7").Value = 0.00684931506849315
$ws.Range("F7").Value = 0.04794520547945205
$ws.Range("J7").Value = 0.1301369863013699
$ws.Range("O7").Value = 0.02054794520547945
$ws.Range("Q7").Value = 0.1986301369863014
$ws.Range("R7").Value = 0.07534246575342465
$ws.Range("S7").Value = 0.4041095890410959
$ws.Range("B8").Value = 0.1130653266331658
$ws.Range("D8").Value = 0.01005025125628141
$ws.Range("E8").Value = 0.002512562814070352
$ws.Range("F8").Value = 0.05778894472361809
$ws.Range("J8").Value = 0.1331658291457286
$ws.Range("O8").Value = 0.02261306532663317
$ws.Range("Q8").Value = 0.1633165829145729
$ws.Range("R8").Value = 0.06532663316582915
$ws.Range("S8").Value = 0.4321608040201005
$ws.Range("B9").Value = 0.09770114942528736
$ws.Range("D9").Value = 0.04022988505747126
$ws.Range("F9").Value = 0.06321839080459771
$ws.Range("J9").Value = 0.1264367816091954
$ws.Range("O9").Value = 0.02298850574712644
$ws.Range("Q9").Value = 0.103448275862069
$ws.Range("R9").Value = 0.06896551724137931
$ws.Range("S9").Value = 0.4770114942528735
$ws.Range("B10").Value = 0.11346316680779
$ws.Range("D10").Value = 0.01608806096528366
$ws.Range("F10").Value = 0.07705334462320068
$ws.Range("J10").Value = 0.1295512277730737
$ws.Range("O10").Value = 0.02201524132091448
$ws.Range("Q10").Value = 0.1761219305673158
$ws.Range("R10").Value = 0.07620660457239628
$ws.Range("S10").Value = 0.3895004233700254
$ws.Range("G11").Value = 0.1752136752136752
$ws.Range("J11").Value = 0.0811965811965812
$ws.Range("K11").Value = 0.2222222222222222
$ws.Range("L11").Value = 0.4957264957264957
$ws.Range("S11").Value = 0.02564102564102564
$ws.Range("G12").Value = 0.71900826446281
$ws.Range("J12").Value = 0.2066115702479339
$ws.Range("L12").Value = 0.03305785123966942
$ws.Range("S12").Value = 0.04132231404958678
$ws.Range("G13").Value = 0.7142857142857143
$ws.Range("J13").Value = 0.25
$ws.Range("S13").Value = 0.03571428571428571
$ws.Range("H15").Value = 0.1666666666666667
$ws.Range("I15").Value = 0.09259259259259259
$ws.Range("J15").Value = 0.3796296296296297
$ws.Range("K15").Value = 0.04166666666666666
$ws.Range("M15").Value = 0.004629629629629629
$ws.Range("O15").Value = 0.05092592592592592
$ws.Range("S15").Value = 0.2638888888888889
$ws.Range("F16").Value = 0.02366863905325444
$ws.Range("H16").Value = 0.1952662721893491
$ws.Range("I16").Value = 0.07100591715976332
$ws.Range("J16").Value = 0.4319526627218935
$ws.Range("K16").Value = 0.1183431952662722
$ws.Range("M16").Value = 0.005917159763313609
$ws.Range("O16").Value = 0.05325443786982249
$ws.Range("S16").Value = 0.1005917159763314
$ws.Range("F17").Value = 0.01169590643274854
$ws.Range("H17").Value = 0.1549707602339181
$ws.Range("I17").Value = 0.0935672514619883
$ws.Range("J17").Value = 0.47953216374269
$ws.Range("K17").Value = 0.08187134502923976
$ws.Range("M17").Value = 0.02631578947368421
$ws.Range("O17").Value = 0.07894736842105263
$ws.Range("S17").Value = 0.07309941520467836
$ws.Range("F18").Value = 0.02515723270440252
$ws.Range("H18").Value = 0.220125786163522
$ws.Range("I18").Value = 0.07547169811320754
$ws.Range("J18").Value = 0.4654088050314465
$ws.Range("K18").Value = 0.03773584905660377
$ws.Range("M18").Value = 0.01257861635220126
$ws.Range("O18").Value = 0.06289308176100629
$ws.Range("S18").Value = 0.1006289308176101
$ws.Range("F19").Value = 0.008748906386701663
$ws.Range("H19").Value = 0.2108486439195101
$ws.Range("I19").Value = 0.08748906386701662
$ws.Range("J19").Value = 0.4111986001749781
$ws.Range("K19").Value = 0.1023622047244094
$ws.Range("M19").Value = 0.01224846894138233
$ws.Range("O19").Value = 0.0804899387576553
$ws.Range("S19").Value = 0.08661417322834646
